$wb = $excel.ActiveWorkbook

# The source sheet was recreated by the automation script (its internal
# sheetId bumps from 1 to 2): add a fresh sheet, drop the old one, and give
# the new sheet back its original name so it lands in the same position.
$oldName = $wb.ActiveSheet.Name
$ws = $wb.Worksheets.Add()
$oldws = $wb.Worksheets.Item($oldName)
$oldws.Delete() | Out-Null
$ws.Name = $oldName

# Header row
$ws.Range("A1").Value = "Content Path"
$ws.Range("B1").Value = "Status"

# Data rows (Status column left blank for these rows)
$ws.Range("A2").Value = "/content/abbvie-pro/de/de/therapy-areas/virology/hepatitis-c/service/thank-you-contact-request"
$ws.Range("A3").Value = "/content/abbvie-pro/de/de/therapy-areas/oncology/hemato-oncology/artikel/closed/direct-downloads/eha23-tag1-cll2"
$ws.Range("A4").Value = "/content/abbvie-pro/de/de/therapy-areas/oncology/hemato-oncology/artikel/closed/direct-downloads/kongresshighlights-ash-2022-t2-x1"
$ws.Range("A5").Value = "/content/abbvie-pro/de/de/therapy-areas/immunology/dermatology/artikel/closed/fobi-kongressbericht-2024"

$ws.Columns.Item(1).ColumnWidth = 117.66666666666667

# Restore the selection to where the user last left it (below the data, at A15)
$ws.Range("A15").Select() | Out-Null
